$d = $word.ActiveDocument
$sec = $d.Sections(1)

function Rename-InlineShape($range, $newName) {
    $inline = $range.InlineShapes(1)
    $shape = $inline.ConvertToShape()
    $shape.Name = $newName
    $shape.ConvertToInlineShape()
}

# Footer "first page" (footer1.xml) - Pearson logo: image1.png -> image2.png
Rename-InlineShape $sec.Footers(2).Range "image2.png"

# Footer "default" (footer2.xml) - Pearson logo: image1.png -> image2.png
Rename-InlineShape $sec.Footers(1).Range "image2.png"

# Header "first page" (header1.xml) - BTec logo: image2.jpg -> image1.jpg
Rename-InlineShape $sec.Headers(2).Range "image1.jpg"
